$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "19"
$ws.Range("D3").Value = "21"
$ws.Range("E3").Value = "3"

$ws.Range("C4").Value = "15"
$ws.Range("D4").Value = "12"
$ws.Range("E4").Value = "2"

$ws.Range("C6").Value = "47"
$ws.Range("D6").Value = "36"
$ws.Range("E6").Value = "4"
$ws.Range("F6").Value = "1"

$ws.Range("C7").Value = "21"
$ws.Range("D7").Value = "13"
$ws.Range("F7").Value = "1"

$ws.Range("C8").Value = "3"
$ws.Range("D8").Value = "5"
$ws.Range("E8").Value = "0"
$ws.Range("F8").Value = "0"

$ws.Range("C9").Value = "1"
$ws.Range("D9").Value = "4"

$ws.Range("C10").Value = "0"
$ws.Range("D10").Value = "2"
$ws.Range("E10").Value = "0"
$ws.Range("F10").Value = "0"

$ws.Range("C11").Value = "29"
$ws.Range("D11").Value = "17"
$ws.Range("F11").Value = "3"

$ws.Range("C13").Value = "28"
$ws.Range("D13").Value = "28"
